# Updated symbol list on Tue Dec 27 10:43:51 UTC 2022 with GitHub Actions
# Refresh the "Price" (column D) and "Volume(1h)" (column E) values on Sheet1
# to the latest scraped coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: Price -> "243.62" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "243.62"
$c.Style = "Normal"

# D4: Price -> "5.386" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "5.386"
$c.Style = "Normal"

# D5: Price -> "0.05957" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "0.05957"
$c.Style = "Normal"

# D6: Price -> "3.431" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "3.431"
$c.Style = "Normal"

# D8: Price -> "0.8109" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.8109"
$c.Style = "Normal"

# D9: Price -> "0.9284" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.9284"
$c.Style = "Normal"

# D11: Price -> "0.07423" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07423"
$c.Style = "Normal"

# D12: Price -> "0.03277" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.03277"
$c.Style = "Normal"

# D13: Price -> "0.03077" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.03077"
$c.Style = "Normal"

# D14: Price -> "0.09359" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "0.09359"
$c.Style = "Normal"

# D15: Price -> "3.848" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "3.848"
$c.Style = "Normal"

# D16: Price -> "0.001586" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.001586"
$c.Style = "Normal"

# D17: Price -> "0.04717" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.04717"
$c.Style = "Normal"

# D18: Price -> "0.0005959" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "0.0005959"
$c.Style = "Normal"

# E18: Volume(1h) label -> "17OneONEWorstin24h"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# D19: Price -> "0.005889" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.005889"
$c.Style = "Normal"

# D20: Price -> "0.001260" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "0.001260"
$c.Style = "Normal"

# E20: Volume(1h) label -> "19BitKanKANBestin24h"
$ws.Range("E20").Value = "19BitKanKANBestin24h"

# D21: Price -> "0.004790" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "0.004790"
$c.Style = "Normal"

# D22: Price -> "0.00007998" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "0.00007998"
$c.Style = "Normal"

# D23: Price -> "3.576" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "3.576"
$c.Style = "Normal"

# D25: Price -> "0.3242" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "0.3242"
$c.Style = "Normal"

# E27: Volume(1h) label -> "26UpBotsUBXT"
$ws.Range("E27").Value = "26UpBotsUBXT"

# D42: Price -> "0.1075" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.1075"
$c.Style = "Normal"

# E43: Volume(1h) label -> "42CEJICEJI"
$ws.Range("E43").Value = "42CEJICEJI"

# D44: Price -> "0.008958" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.008958"
$c.Style = "Normal"

# D47: Price -> "0.6849" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.6849"
$c.Style = "Normal"

# D48: Price -> "0.002069" (keep as text, matching the sheet's inline-string cells)
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.002069"
$c.Style = "Normal"

